$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row value updates (Price column D, Volume(1h) column E).
# A handful of D-column prices are plain decimals (e.g. "559.73") that Excel
# would otherwise auto-convert to a number; prefix them with an apostrophe
# (quote-prefix) so they are entered as text, matching the source data which
# stores every Price/Volume cell as a text string.

$ws.Range("D2").Value = "68.876.12"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.441.31"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'559.73"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'162.31"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  +8.46%  "
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").Value = "'0.332"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("E12").Value = "  -4.77%  "
$ws.Range("E13").Value = "  +4.78%  "
$ws.Range("D14").Value = "68.760.46"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "2.889.53"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "'23.31"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "2.441.77"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "'10.62"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "'339.11"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'6.98"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +3.42%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'66.13"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").Value = "'3.74"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").Value = "2.567.51"
$ws.Range("E26").Value = "  -1.16%  "

# Rows 27 and 28 swap places (Aptos <-> Binance-PegBSC-USD) with updated
# price/volume; the rank index in column A stays tied to the row, not the coin.
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.02"
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'8.25"
$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").Value = "'7.16"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'430.23"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").Value = "'160.42"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "'0.105"
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("D40").Value = "'0.301"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").Value = "'1.51"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").Value = "'2.08"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("D46").Value = "'130.78"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "'0.0721"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("E50").Value = "  +2.96%  "
$ws.Range("D51").Value = "'0.0923"
$ws.Range("E51").Value = "  +0.58%  "
